$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 985.8333
$ws.Range("J17").Value = 999.46155
$ws.Range("L17").Value = 2998.38465
$ws.Range("N17").Value = -3334.38465
$ws.Range("H55").Value = 318.05405
$ws.Range("J55").Value = 346.5
$ws.Range("L55").Value = 346.5
$ws.Range("N55").Value = -774.5
$ws.Range("H95").Value = 39641.332
$ws.Range("J95").Value = 39641.332
$ws.Range("L95").Value = 39641.332
$ws.Range("N95").Value = -45133.332
$ws.Range("H112").Value = 2926.0557
$ws.Range("I112").Value = 1020
$ws.Range("J112").Value = 3233.484
$ws.Range("K112").Value = 3060
$ws.Range("L112").Value = 9700.451999999999
$ws.Range("M112").Value = -1952
$ws.Range("N112").Value = -11916.452
$ws.Range("H129").Value = 1233.1163
$ws.Range("I129").Value = 501.375
$ws.Range("J129").Value = 1666.7407
$ws.Range("K129").Value = 1504.125
$ws.Range("L129").Value = 5000.2221
$ws.Range("M129").Value = 3495.875
$ws.Range("N129").Value = -15000.2221
$ws.Range("H132").Value = 4478.9736
$ws.Range("I132").Value = 3647.4285
$ws.Range("J132").Value = 6807.3
$ws.Range("K132").Value = 10942.2855
$ws.Range("L132").Value = 20421.9
$ws.Range("M132").Value = -8412.2855
$ws.Range("N132").Value = -25481.9
$ws.Range("H138").Value = 1156.2699
$ws.Range("I138").Value = 924.4211
$ws.Range("J138").Value = 1508.68
$ws.Range("K138").Value = 2773.2633
$ws.Range("L138").Value = 4526.04
$ws.Range("M138").Value = 2366.7367
$ws.Range("N138").Value = -14806.04
$ws.Range("H141").Value = 2892.8513
$ws.Range("I141").Value = 734.6462
$ws.Range("J141").Value = 18479.889
$ws.Range("K141").Value = 2203.9386
$ws.Range("L141").Value = 55439.667
$ws.Range("M141").Value = 2976.0614
$ws.Range("N141").Value = -65799.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1139.17
$ws.Range("I32").Value = 1011.7765
$ws.Range("J32").Value = 1861.0667
$ws.Range("K32").Value = 1011.7765
$ws.Range("L32").Value = 1861.0667
$ws.Range("M32").Value = -724.7765000000001
$ws.Range("N32").Value = -2435.0667
$ws.Range("H61").Value = 1810.4762
$ws.Range("I61").Value = 1865.3334
$ws.Range("J61").Value = 1769.3334
$ws.Range("K61").Value = 1865.3334
$ws.Range("L61").Value = 1769.3334
$ws.Range("M61").Value = -1653.3334
$ws.Range("N61").Value = -2193.3334
$ws.Range("H74").Value = 1060.638
$ws.Range("I74").Value = 847.8158
$ws.Range("K74").Value = 847.8158
$ws.Range("M74").Value = 26.18420000000003
$ws.Range("H77").Value = 1060.638
$ws.Range("I77").Value = 847.8158
$ws.Range("K77").Value = 4239.079
$ws.Range("M77").Value = 128.9210000000003
$ws.Range("H132").Value = 3716.5
$ws.Range("I132").Value = 2360.6667
$ws.Range("K132").Value = 7082.000100000001
$ws.Range("M132").Value = -4552.000100000001
$ws.Range("H136").Value = 1810.4762
$ws.Range("I136").Value = 1865.3334
$ws.Range("J136").Value = 1769.3334
$ws.Range("K136").Value = 5596.0002
$ws.Range("L136").Value = 5308.0002
$ws.Range("M136").Value = -3046.0002
$ws.Range("N136").Value = -10408.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4000.2205
$ws.Range("I134").Value = 1491.6451
$ws.Range("J134").Value = 6777.5713
$ws.Range("K134").Value = 4474.9353
$ws.Range("L134").Value = 20332.7139
$ws.Range("M134").Value = -1939.9353
$ws.Range("N134").Value = -25402.7139

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5651780.5
$ws.Range("I31").Value = 1587.8043
$ws.Range("J31").Value = 25644770
$ws.Range("K31").Value = 1587.8043
$ws.Range("L31").Value = 25644770
$ws.Range("M31").Value = -1292.8043
$ws.Range("N31").Value = -25645360
$ws.Range("H34").Value = 5651780.5
$ws.Range("I34").Value = 1587.8043
$ws.Range("J34").Value = 25644770
$ws.Range("K34").Value = 1587.8043
$ws.Range("L34").Value = 25644770
$ws.Range("M34").Value = -1385.8043
$ws.Range("N34").Value = -25645174
$ws.Range("H132").Value = 1872.5193
$ws.Range("I132").Value = 951.9
$ws.Range("J132").Value = 3127.9092
$ws.Range("K132").Value = 2855.7
$ws.Range("L132").Value = 9383.7276
$ws.Range("M132").Value = -325.6999999999998
$ws.Range("N132").Value = -14443.7276
$ws.Range("H134").Value = 1494.4728
$ws.Range("I134").Value = 621.6429000000001
$ws.Range("J134").Value = 2399.6296
$ws.Range("K134").Value = 1864.9287
$ws.Range("L134").Value = 7198.888800000001
$ws.Range("M134").Value = 670.0712999999998
$ws.Range("N134").Value = -12268.8888

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 512.0833
$ws.Range("I5").Value = 232.42105
$ws.Range("J5").Value = 1574.8
$ws.Range("K5").Value = 697.26315
$ws.Range("L5").Value = 4724.4
$ws.Range("M5").Value = -585.26315
$ws.Range("N5").Value = -4948.4
$ws.Range("H8").Value = 86
$ws.Range("I8").Value = 86
$ws.Range("K8").Value = 258
$ws.Range("M8").Value = -119
$ws.Range("H99").Value = 1700
$ws.Range("I99").Value = 2000
$ws.Range("J99").Value = 1625
$ws.Range("K99").Value = 6000
$ws.Range("L99").Value = 4875
$ws.Range("M99").Value = -3754
$ws.Range("N99").Value = -9367
$ws.Range("H135").Value = 512.0833
$ws.Range("I135").Value = 232.42105
$ws.Range("J135").Value = 1574.8
$ws.Range("K135").Value = 2091.78945
$ws.Range("L135").Value = 14173.2
$ws.Range("M135").Value = 443.2105499999998
$ws.Range("N135").Value = -19243.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1070322
$ws.Range("I132").Value = 2316078.2
$ws.Range("J132").Value = 2530.9524
$ws.Range("K132").Value = 6948234.600000001
$ws.Range("L132").Value = 7592.8572
$ws.Range("M132").Value = -6945704.600000001
$ws.Range("N132").Value = -12652.8572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 15168504
$ws.Range("I132").Value = 26343836
$ws.Range("J132").Value = 1983.5714
$ws.Range("K132").Value = 79031508
$ws.Range("L132").Value = 5950.7142
$ws.Range("M132").Value = -79028978
$ws.Range("N132").Value = -11010.7142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1384.7931
$ws.Range("I132").Value = 1288.9117
$ws.Range("J132").Value = 1520.625
$ws.Range("K132").Value = 3866.7351
$ws.Range("L132").Value = 4561.875
$ws.Range("M132").Value = -1336.7351
$ws.Range("N132").Value = -9621.875
